$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" and "is_enabled" dictionary-field columns (D:E) from
# the generated import-template header row. Deleting the entire columns
# shifts the following "order_by" / "rem" columns left so the row again ends
# at column E (F:G no longer exist), matching the trimmed column set.
$ws.Range("D1:E1").EntireColumn.Delete()
